$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplemental transect draw post power analysis: two stations were
# renumbered in the redraw -- "LTI" became "LTI3" and "LTI10" became
# "LTI12" (and their paired TRAN_ST/TRAN_END LABEL values likewise).
$ws.Range("B33").Value = "LTI3"
$ws.Range("J33").Value = "LTI3-2"
$ws.Range("B34").Value = "LTI3"
$ws.Range("J34").Value = "LTI3-2"

$ws.Range("B35").Value = "LTI12"
$ws.Range("J35").Value = "LTI12-1"
$ws.Range("B36").Value = "LTI12"
$ws.Range("J36").Value = "LTI12-1"

# Leave the view scrolled/selected where the editor ended up working.
$excel.Goto($ws.Range("K36"), $true)
